$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply D82:D87 as one shared formula (matches the row above's shared-formula pattern)
$ws.Range("D82:D87").Formula = "=C82/(24*60)"

# New row 88: 四方坪站 (station 2) for date 45944
$ws.Range("A88").Value = 45944
$ws.Range("B88").Value = "四方坪站"
$ws.Range("C88").Formula = "=15417/126"
$ws.Range("D88").Formula = "=C88/(24*60)"
$ws.Range("E88").Formula = "=8775.35/126"
$ws.Range("F88").Formula = "=3053.72/126"
$ws.Range("G88").Formula = "=8775.35/(15417/60)"
$ws.Range("H88").Formula = "=378/126"

# New row 89: 高岭站 (station 3) for date 45944
$ws.Range("A89").Value = 45944
$ws.Range("B89").Value = "高岭站"
$ws.Range("C89").Formula = "=6107/36"
$ws.Range("D89").Formula = "=C89/(24*60)"
$ws.Range("E89").Formula = "=4056.73/36"
$ws.Range("F89").Formula = "=1081.36/36"
$ws.Range("G89").Formula = "=4056.73/(6107/60)"
$ws.Range("H89").Formula = "=160/36"

# Update selection to match the new active cell
[void]$ws.Range("I89").Select()
